$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table data (player, position, team) replacing the old roster listing.
$data = @(
    ,("Russell Westbrook", "PG,SG", "Denver Nuggets")
    ,("Ayo Dosunmu", "PG,SG,SF", "Chicago Bulls")
    ,("Payton Pritchard", "PG,SG", "Boston Celtics")
    ,("Nick Smith Jr.", "PG,SG", "Charlotte Hornets")
    ,("Jaylen Brown", "SG,SF", "Boston Celtics")
    ,("Paolo Banchero", "SF,PF", "Orlando Magic")
    ,("Deni Avdija", "SF,PF", "Portland Trail Blazers")
    ,("Nikola Jokic", "C", "Denver Nuggets")
    ,("Rudy Gobert", "C", "Minnesota Timberwolves")
    ,("Dejounte Murray", "PG,SG", "New Orleans Pelicans")
    ,("Jose Alvarado", "PG", "New Orleans Pelicans")
    ,("Jakob Poeltl", "C", "Toronto Raptors")
    ,("Jalen Green", "PG,SG", "Houston Rockets")
    ,("Chris Paul", "PG", "San Antonio Spurs")
    ,("Pascal Siakam", "SF,PF,C", "Indiana Pacers")
    ,("Chet Holmgren", "PF,C", "Oklahoma City Thunder")
    ,("Jalen Suggs", "PG,SG", "Orlando Magic")
    ,("Khris Middleton", "SF", "Milwaukee Bucks")
)

$rowCount = $data.Count

# Clear out any existing data rows below the header first (old table had 17 rows of data).
$oldDataRange = $ws.Range("A2:C18")
$oldDataRange.ClearContents()

# Write the new data starting at row 2.
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
